$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Range("B2").Value = "test"
Write-Host "done"
